# "Generate Report for Archive"
# - Update localization status text from "Ready for handoff" to "In Translation"
#   (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all shared that string).
# - Narrow the per-language "Status" columns (Overview E:F, zh-cn C, de-de C)
#   from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text update: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Column width update for the Status columns ---
# Target stored width is 13.4101845877511 characters; ColumnWidth assignments
# get snapped to the host's display-pixel grid, so back off by the fixed
# 5/6-character padding before assigning to land on the closest pixel.
$newStatusWidth = 13.4101845877511 - (5 / 6)

$wsOverview.Range("E1:F1").ColumnWidth = $newStatusWidth
$wsZhCn.Range("C1").ColumnWidth = $newStatusWidth
$wsDeDe.Range("C1").ColumnWidth = $newStatusWidth
